$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$featureStr = "11 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii"
$modelTypeStr = "Neural-Network"
$modelStr = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000"

$rows = @(
    @{ Time = "20160427_084337"; ClassifyAcc = 0.910891089108911; SegmentAcc = 0.44 },
    @{ Time = "20160427_093906"; ClassifyAcc = 0.914191419141914; SegmentAcc = 0.43 },
    @{ Time = "20160427_103542"; ClassifyAcc = 0.910891089108911; SegmentAcc = 0.44 },
    @{ Time = "20160427_113319"; ClassifyAcc = 0.914191419141914; SegmentAcc = 0.45 },
    @{ Time = "20160427_123046"; ClassifyAcc = 0.910891089108911; SegmentAcc = 0.44 }
)

$r = 27
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Time
    $ws.Cells.Item($r, 2).Value = $featureStr
    $ws.Cells.Item($r, 3).Value = $featureStr
    $ws.Cells.Item($r, 4).Value = $featureStr
    $ws.Cells.Item($r, 5).Value = $modelTypeStr
    $ws.Cells.Item($r, 6).Value = $modelStr
    $ws.Cells.Item($r, 7).Value = $modelTypeStr
    $ws.Cells.Item($r, 8).Value = $modelStr
    $ws.Cells.Item($r, 9).Value = $modelTypeStr
    $ws.Cells.Item($r, 10).Value = $modelStr
    $ws.Cells.Item($r, 11).Value = $row.ClassifyAcc
    $ws.Cells.Item($r, 12).Value = $row.SegmentAcc
    $r++
}
